# Auto-update predictions and index for 2025-10-15
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as text (avoids Excel auto-converting numeric-looking
# strings like "100" or "1.40" into real numbers), then restore the cell's
# default style so no stray number-format style is left behind.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Update the fixtures index count on row 3
$ws.Range("E3").Value = "28 of 28 Tips"

# Rows 9-12 shift: a new fixture (Palmeiras vs Red Bull Bragantino) is
# inserted at row 9, pushing the following fixtures down by one row, and a
# new fixture (USC Landhaus Women vs Slavia Prague Women) lands in row 12.

# Row 9: Palmeiras vs Red Bull Bragantino
$ws.Range("A9").Value = "Palmeiras vs Red Bull Bragantino"
$ws.Range("B9").Value = "Palmeiras"
$ws.Range("C9").Value = "Brazil"
$ws.Range("D9").Value = "15th Oct 23:00"
$ws.Range("E9").Value = "19 of 19 Tips"
Set-TextValue "F9" "100"
Set-TextValue "G9" "1.40"

# Row 10: Sevilla vs Mallorca
$ws.Range("A10").Value = "Sevilla vs Mallorca"
$ws.Range("B10").Value = "Sevilla"
$ws.Range("C10").Value = "Spain"
$ws.Range("D10").Value = "18th Oct 13:00"
$ws.Range("E10").Value = "14 of 17 Tips"
Set-TextValue "F10" "82"
Set-TextValue "G10" "1.91"

# Row 11: Chelsea Women vs Paris FC Women
$ws.Range("A11").Value = "Chelsea Women vs Paris FC Women"
$ws.Range("B11").Value = "Chelsea Women"
$ws.Range("C11").Value = "Europe"
$ws.Range("D11").Value = "15th Oct 20:00"
$ws.Range("E11").Value = "13 of 14 Tips"
Set-TextValue "F11" "93"
Set-TextValue "G11" "1.14"

# Row 12: USC Landhaus Women vs Slavia Prague Women
$ws.Range("A12").Value = "USC Landhaus Women vs Slavia Prague Women"
$ws.Range("B12").Value = "USC Landhaus Women"
$ws.Range("C12").Value = "Europe"
$ws.Range("D12").Value = "15th Oct 18:45"
$ws.Range("E12").Value = "13 of 14 Tips"
Set-TextValue "F12" "93"
Set-TextValue "G12" "1.90"
